$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely (naive component forecaster bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Update forecast values with corrected precision
$ws.Range("E3").Value = 5.080273296954396
$ws.Range("E4").Value = 2.807231216534301
$ws.Range("C5").Value = -3.942037578692492
$ws.Range("E5").Value = -1.648748515828491
$ws.Range("C6").Value = -0.9140166223623458
$ws.Range("E6").Value = 1.821983295885099
$ws.Range("C7").Value = -2.839753013810498
$ws.Range("E7").Value = -1.632723506456923
$ws.Range("C9").Value = 4.960109259035428
$ws.Range("C13").Value = -2.700325749999499
$ws.Range("E13").Value = -0.3858735870725494
$ws.Range("C14").Value = -1.479696720105184
$ws.Range("C15").Value = 5.469647210234974
$ws.Range("E15").Value = 3.061326532789521
$ws.Range("E16").Value = -0.6155071485167807
$ws.Range("C19").Value = 3.458696398997096
$ws.Range("E19").Value = 2.610227683091315
$ws.Range("C20").Value = 2.192778679161966
$ws.Range("E20").Value = -0.5835597102573087
$ws.Range("C21").Value = 2.77241330895972
$ws.Range("C22").Value = 3.408364488606752
$ws.Range("C23").Value = 1.376958470962353
$ws.Range("C24").Value = 3.14581984265847
$ws.Range("C25").Value = 2.594480907596486
$ws.Range("E25").Value = 2.97885695339557
$ws.Range("E26").Value = 3.056075254340018
$ws.Range("C27").Value = 3.665688413913704
$ws.Range("E27").Value = 3.126710782028064
$ws.Range("E28").Value = 4.124307769579483
$ws.Range("C29").Value = 4.143226503463815
$ws.Range("E32").Value = 4.888255652935958
$ws.Range("C34").Value = 1.666553973046025
$ws.Range("E34").Value = -1.376301649685407
$ws.Range("C35").Value = 0.399474938574329
$ws.Range("E35").Value = 1.894543511868685
$ws.Range("E36").Value = 2.95288809451808
$ws.Range("C37").Value = 1.61478104109658
$ws.Range("E37").Value = 1.61617062332684
$ws.Range("C38").Value = 1.879266440112781
$ws.Range("E38").Value = -0.5015683214423916
$ws.Range("C39").Value = -1.647049671756318
$ws.Range("C41").Value = -2.367053688984511
$ws.Range("C42").Value = -2.620683231370935
$ws.Range("E42").Value = -3.531225750971467
$ws.Range("C43").Value = -4.823846530890474
$ws.Range("C45").Value = -2.690210378056657
$ws.Range("E45").Value = -0.5080092691694071
$ws.Range("C46").Value = -3.036556262700263
$ws.Range("C47").Value = -3.632744220111894
$ws.Range("E47").Value = -3.88801767603163
$ws.Range("C48").Value = -0.5744163079740128
$ws.Range("E48").Value = -0.6956477387308979
$ws.Range("C49").Value = -2.862797582711463
$ws.Range("E49").Value = -2.931248208044468
$ws.Range("E51").Value = -2.64937131332158
$ws.Range("C52").Value = -0.1892239049850142
$ws.Range("C53").Value = -1.201318493118475
$ws.Range("E53").Value = -2.375407956079478
